$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: change from duplicate "богдана" row to "Підсумок" summary row (values from old row4, but B=1)
$ws.Range("A3").Value = "Підсумок"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0

# Row 4: becomes the label row (previously row 5); clear leftover cells from old row 4 (E4:L4)
$ws.Range("A4").Value = "Будуть"
$ws.Range("B4").Value = "Не будуть"
$ws.Range("C4").Value = "Не знаю"
$ws.Range("D4").Value = "Відмітилось"
$ws.Range("E4:L4").ClearContents()

# Row 5: becomes the counts row (previously row 6), with updated values
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 1

# Remove the old row 6 entirely (shift rows up, deleting the last row)
$ws.Rows.Item(6).Delete()
